$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.329.15"
$ws.Range("E2").Value = "'  -4.90%  "
$ws.Range("D3").Value = "'1.559.94"
$ws.Range("E3").Value = "'  -5.22%  "
$ws.Range("E4").Value = "'  +0.18%  "
$ws.Range("E5").Value = "'  +0.14%  "
$ws.Range("D6").Value = "'290.00"
$ws.Range("E6").Value = "'  -3.54%  "
$ws.Range("E7").Value = "'  -2.08%  "
$ws.Range("D8").Value = "'49.13"
$ws.Range("E8").Value = "'  -2.90%  "
$ws.Range("D9").Value = "'0.3406"
$ws.Range("E9").Value = "'  -2.74%  "
$ws.Range("D10").Value = "'1.162"
$ws.Range("E10").Value = "'  -4.42%  "
$ws.Range("D11").Value = "'0.07630"
$ws.Range("E11").Value = "'  -5.40%  "
$ws.Range("E12").Value = "'  +0.20%  "
$ws.Range("D13").Value = "'21.39"
$ws.Range("E13").Value = "'  -3.17%  "
$ws.Range("D14").Value = "'6.032"
$ws.Range("E14").Value = "'  -4.34%  "
$ws.Range("D15").Value = "'6.904"
$ws.Range("E15").Value = "'  -4.73%  "
$ws.Range("D16").Value = "'1.562.10"
$ws.Range("E16").Value = "'  -4.89%  "
$ws.Range("D17").Value = "'0.00001125"
$ws.Range("E17").Value = "'  -7.18%  "
$ws.Range("D18").Value = "'89.88"
$ws.Range("E18").Value = "'  -5.42%  "
$ws.Range("D19").Value = "'0.06715"
$ws.Range("E19").Value = "'  -3.62%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "'  +0.12%  "
$ws.Range("D21").Value = "'6.216"
$ws.Range("E21").Value = "'  -6.13%  "
$ws.Range("D22").Value = "'16.50"
$ws.Range("E22").Value = "'  -5.31%  "
$ws.Range("D23").Value = "'0.5293"
$ws.Range("E23").Value = "'  -7.49%  "
$ws.Range("D24").Value = "'11.96"
$ws.Range("E24").Value = "'  -4.01%  "
$ws.Range("D25").Value = "'22.325.85"
$ws.Range("E25").Value = "'  -4.95%  "
$ws.Range("D26").Value = "'2.394"
$ws.Range("E26").Value = "'  -1.08%  "
$ws.Range("D27").Value = "'2.817"
$ws.Range("E27").Value = "'  -5.14%  "
$ws.Range("D28").Value = "'20.17"
$ws.Range("E28").Value = "'  -4.04%  "
$ws.Range("D29").Value = "'145.90"
$ws.Range("D30").Value = "'4.983"
$ws.Range("E30").Value = "'  -3.87%  "
$ws.Range("D31").Value = "'125.20"
$ws.Range("E31").Value = "'  -4.93%  "
$ws.Range("D32").Value = "'1.733.77"
$ws.Range("E32").Value = "'  -5.05%  "
$ws.Range("D33").Value = "'6.155"
$ws.Range("E33").Value = "'  -9.95%  "
$ws.Range("B34").Value = "'ImmutableX"
$ws.Range("C34").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.005"
$ws.Range("E34").Value = "'  +1.64%  "
$ws.Range("B35").Value = "'WEMIXTOKEN"
$ws.Range("C35").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.009"
$ws.Range("E35").Value = "'  -6.06%  "
$ws.Range("D36").Value = "'10.03"
$ws.Range("E36").Value = "'  -10.24%  "
$ws.Range("D37").Value = "'0.08475"
$ws.Range("E37").Value = "'  -3.57%  "
$ws.Range("D38").Value = "'0.02535"
$ws.Range("E38").Value = "'  -5.86%  "
$ws.Range("D39").Value = "'0.2310"
$ws.Range("E39").Value = "'  -4.51%  "
$ws.Range("D40").Value = "'5.508"
$ws.Range("E40").Value = "'  -6.81%  "
$ws.Range("B41").Value = "'TrustWalletToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.320"
$ws.Range("E41").Value = "'  +1.94%  "
$ws.Range("B42").Value = "'Hedera"
$ws.Range("C42").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.06371"
$ws.Range("E42").Value = "'  -5.98%  "
$ws.Range("D43").Value = "'11.66"
$ws.Range("E43").Value = "'  -9.15%  "
$ws.Range("D44").Value = "'0.6330"
$ws.Range("E44").Value = "'  -7.85%  "
$ws.Range("B45").Value = "'Frax"
$ws.Range("C45").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "'  +0.12%  "
$ws.Range("B46").Value = "'EnergySwap"
$ws.Range("C46").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'14.06"
$ws.Range("E46").Value = "'  -9.74%  "
$ws.Range("D47").Value = "'0.5963"
$ws.Range("E47").Value = "'  -6.52%  "
$ws.Range("D48").Value = "'3.755"
$ws.Range("E48").Value = "'  -4.43%  "
$ws.Range("D49").Value = "'2.087"
$ws.Range("E49").Value = "'  -7.01%  "
$ws.Range("D50").Value = "'1.261"
$ws.Range("E50").Value = "'  +2.45%  "
$ws.Range("D51").Value = "'123.80"
$ws.Range("E51").Value = "'  -2.59%  "
